$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (rows 2-21) is cleared entirely in the target state
$ws.Range("B2:B21").ClearContents()

# Columns A and C are updated with the new ranking values
$ws.Range("A2").Value = "华胜天成"
$ws.Range("C2").Value = "华胜天成"
$ws.Range("A3").Value = "兆易创新"
$ws.Range("C3").Value = "天奇股份"
$ws.Range("A4").Value = "天奇股份"
$ws.Range("C4").Value = "博纳影业"
$ws.Range("A5").Value = "光线传媒"
$ws.Range("C5").Value = "巨力索具"
$ws.Range("A6").Value = "捷成股份"
$ws.Range("C6").Value = "利欧股份"
$ws.Range("A7").Value = "国安股份"
$ws.Range("C7").Value = "光线传媒"
$ws.Range("A8").Value = "利欧股份"
$ws.Range("C8").Value = "协鑫集成"
$ws.Range("A9").Value = "深科技"
$ws.Range("C9").Value = "嘉美包装"
$ws.Range("A10").Value = "博纳影业"
$ws.Range("C10").Value = "大位科技"
$ws.Range("A11").Value = "紫金矿业"
$ws.Range("C11").Value = "掌阅科技"
$ws.Range("A12").Value = "掌阅科技"
$ws.Range("C12").Value = "汉缆股份"
$ws.Range("A13").Value = "汉缆股份"
$ws.Range("C13").Value = "兆易创新"
$ws.Range("A14").Value = "巨力索具"
$ws.Range("C14").Value = "二六三"
$ws.Range("A15").Value = "澜起科技"
$ws.Range("C15").Value = "紫金矿业"
$ws.Range("A16").Value = "协鑫集成"
$ws.Range("C16").Value = "万向钱潮"
$ws.Range("A17").Value = "洛阳钼业"
$ws.Range("C17").Value = "风语筑"
$ws.Range("A18").Value = "贵州茅台"
$ws.Range("C18").Value = "深科技"
$ws.Range("A19").Value = "嘉美包装"
$ws.Range("C19").Value = "洛阳钼业"
$ws.Range("A20").Value = "东方财富"
$ws.Range("C20").Value = "捷成股份"
$ws.Range("A21").Value = "风语筑"
$ws.Range("C21").Value = "国安股份"
